# Update Sheets via scheduled runner: refresh market price / profit values
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1357.75
$ws.Range("I28").Value = 723.8571
$ws.Range("J28").Value = 5795
$ws.Range("K28").Value = 723.8571
$ws.Range("L28").Value = 5795
$ws.Range("M28").Value = -238.8570999999999
$ws.Range("N28").Value = -6765
# Row 106
$ws.Range("H106").Value = 38462430
$ws.Range("I106").Value = 45455080
$ws.Range("K106").Value = 45455080
$ws.Range("M106").Value = -45454449
# Row 112
$ws.Range("H112").Value = 5636.967
$ws.Range("J112").Value = 5932.4644
$ws.Range("L112").Value = 17797.3932
$ws.Range("N112").Value = -20013.3932
# Row 137
$ws.Range("H137").Value = 107775.88
$ws.Range("I137").Value = 138861.61
$ws.Range("J137").Value = 6747.25
$ws.Range("K137").Value = 416584.83
$ws.Range("L137").Value = 20241.75
$ws.Range("M137").Value = -414034.83
$ws.Range("N137").Value = -25341.75
# Row 138
$ws.Range("H138").Value = 5013.0244
$ws.Range("I138").Value = 2909.111
$ws.Range("J138").Value = 5604.75
$ws.Range("K138").Value = 8727.332999999999
$ws.Range("L138").Value = 16814.25
$ws.Range("M138").Value = -3587.332999999999
$ws.Range("N138").Value = -27094.25
# Row 141
$ws.Range("H141").Value = 17004.285
$ws.Range("I141").Value = 27011.75
$ws.Range("K141").Value = 81035.25
$ws.Range("M141").Value = -75855.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9271.697
$ws.Range("I32").Value = 8311.6045
$ws.Range("K32").Value = 8311.6045
$ws.Range("M32").Value = -8024.604499999999
# Row 44
$ws.Range("H44").Value = 8161.5
$ws.Range("J44").Value = 8161.5
$ws.Range("L44").Value = 8161.5
$ws.Range("N44").Value = -9137.5
# Row 74
$ws.Range("H74").Value = 73703.16
$ws.Range("I74").Value = 4805.8887
$ws.Range("K74").Value = 4805.8887
$ws.Range("M74").Value = -3931.8887
# Row 77
$ws.Range("H77").Value = 73703.16
$ws.Range("I77").Value = 4805.8887
$ws.Range("K77").Value = 24029.4435
$ws.Range("M77").Value = -19661.4435

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 16670557
$ws.Range("I86").Value = 24077060
$ws.Range("K86").Value = 24077060
$ws.Range("M86").Value = -24075937
# Row 89
$ws.Range("H89").Value = 16670557
$ws.Range("I89").Value = 24077060
$ws.Range("K89").Value = 120385300
$ws.Range("M89").Value = -120379684
# Row 134
$ws.Range("H134").Value = 11013.263
$ws.Range("I134").Value = 2959.1
$ws.Range("J134").Value = 19962.334
$ws.Range("K134").Value = 8877.3
$ws.Range("L134").Value = 59887.00199999999
$ws.Range("M134").Value = -6342.299999999999
$ws.Range("N134").Value = -64957.00199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 1004.5
$ws.Range("I14").Value = 1004.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1004.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -834.5
$ws.Range("N14").ClearContents()
# Row 16
$ws.Range("H16").Value = 652.6667
$ws.Range("I16").Value = 406
$ws.Range("K16").Value = 406
$ws.Range("M16").Value = -119
# Row 22
$ws.Range("H22").Value = 142859300
$ws.Range("I22").Value = 2919
$ws.Range("J22").Value = 500000260
$ws.Range("K22").Value = 2919
$ws.Range("L22").Value = 500000260
$ws.Range("M22").Value = -2569
$ws.Range("N22").Value = -500000960
# Row 31
$ws.Range("H31").Value = 39130.15
$ws.Range("I31").Value = 3701.75
$ws.Range("J31").Value = 67472.87
$ws.Range("K31").Value = 3701.75
$ws.Range("L31").Value = 67472.87
$ws.Range("M31").Value = -3406.75
$ws.Range("N31").Value = -68062.87
# Row 34
$ws.Range("H34").Value = 39130.15
$ws.Range("I34").Value = 3701.75
$ws.Range("J34").Value = 67472.87
$ws.Range("K34").Value = 3701.75
$ws.Range("L34").Value = 67472.87
$ws.Range("M34").Value = -3499.75
$ws.Range("N34").Value = -67876.87
# Row 99
$ws.Range("H99").Value = 4821.6665
$ws.Range("I99").Value = 4233.3335
$ws.Range("J99").Value = 5998.3335
$ws.Range("K99").Value = 4233.3335
$ws.Range("L99").Value = 5998.3335
$ws.Range("M99").Value = -2735.3335
$ws.Range("N99").Value = -8994.3335
# Row 105
$ws.Range("H105").Value = 874.7059
$ws.Range("I105").Value = 858
$ws.Range("K105").Value = 858
$ws.Range("M105").Value = 889
# Row 113
$ws.Range("H113").Value = 652.6667
$ws.Range("I113").Value = 406
$ws.Range("K113").Value = 406
$ws.Range("M113").Value = 1764
# Row 126
$ws.Range("H126").Value = 4821.6665
$ws.Range("I126").Value = 4233.3335
$ws.Range("J126").Value = 5998.3335
$ws.Range("K126").Value = 12700.0005
$ws.Range("L126").Value = 17995.0005
$ws.Range("M126").Value = -10230.0005
$ws.Range("N126").Value = -22935.0005
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 1883.46
$ws.Range("I134").Value = 1209.2433
$ws.Range("K134").Value = 3627.7299
$ws.Range("M134").Value = -1092.7299
# Row 141
$ws.Range("H141").Value = 202812.81
$ws.Range("J141").Value = 202812.81
$ws.Range("L141").Value = 202812.81
$ws.Range("N141").Value = -213172.81

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 150
$ws.Range("I39").Value = 150
$ws.Range("K39").Value = 450
$ws.Range("M39").Value = -156
# Row 113
$ws.Range("H113").Value = 2991
$ws.Range("I113").Value = 2983.3333
$ws.Range("K113").Value = 8949.999899999999
$ws.Range("M113").Value = -6779.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 4933.3335
$ws.Range("I20").Value = 4800
$ws.Range("K20").Value = 4800
$ws.Range("M20").Value = -4555
# Row 24
$ws.Range("H24").Value = 21666.666
$ws.Range("I24").Value = 9000
$ws.Range("K24").Value = 9000
$ws.Range("M24").Value = -8827
# Row 102
$ws.Range("H102").Value = 5398901
$ws.Range("I102").Value = 6945564
$ws.Range("K102").Value = 6945564
$ws.Range("M102").Value = -6943942
# Row 107
$ws.Range("H107").Value = 7589.357
$ws.Range("I107").Value = 12994
$ws.Range("J107").Value = 383.16666
$ws.Range("K107").Value = 12994
$ws.Range("L107").Value = 383.16666
$ws.Range("M107").Value = -11074
$ws.Range("N107").Value = -4223.16666
# Row 113
$ws.Range("H113").Value = 8385752
$ws.Range("J113").Value = 3170.182
$ws.Range("L113").Value = 3170.182
$ws.Range("N113").Value = -7510.182

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 65225.145
$ws.Range("J22").Value = 3375
$ws.Range("L22").Value = 3375
$ws.Range("N22").Value = -3965
# Row 27
$ws.Range("H27").Value = 65225.145
$ws.Range("J27").Value = 3375
$ws.Range("L27").Value = 3375
$ws.Range("N27").Value = -3589
# Row 55
$ws.Range("H55").Value = 1151.3334
$ws.Range("I55").Value = 958.55554
$ws.Range("J55").Value = 1440.5
$ws.Range("K55").Value = 958.55554
$ws.Range("L55").Value = 1440.5
$ws.Range("M55").Value = -785.55554
$ws.Range("N55").Value = -1786.5
# Row 136
$ws.Range("H136").Value = 147496.14
$ws.Range("I136").Value = 170787.17
$ws.Range("K136").Value = 512361.51
$ws.Range("M136").Value = -509811.51

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 25999.666
$ws.Range("I31").Value = 35999
$ws.Range("K31").Value = 35999
$ws.Range("M31").Value = -35651
# Row 107
$ws.Range("H107").Value = 30304102
$ws.Range("I107").Value = 40000540
$ws.Range("K107").Value = 120001620
$ws.Range("M107").Value = -119999700
# Row 132
$ws.Range("H132").Value = 23514178
$ws.Range("I132").Value = 27783052
$ws.Range("K132").Value = 83349156
$ws.Range("M132").Value = -83346626

